# Apply underline formatting to the five "sub-question" list paragraphs
# in the research document, plus the empty paragraph immediately
# preceding the "maintainability" question (whose paragraph mark also
# picked up the underline run property in the source edit).

$d = $word.ActiveDocument

function Underline-Paragraph($needleText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $needleText"
        return $null
    }
    $para = $rng.Paragraphs(1)
    $para.Range.Font.Underline = 1
    return $para
}

function Underline-EmptyParagraph($para) {
    # Formatting an empty (run-less) paragraph range directly does not
    # stick, so nudge it: insert a throwaway character, format the
    # resulting range (character + paragraph mark), then delete just
    # the throwaway character again, leaving the paragraph mark's
    # run properties (<w:pPr><w:rPr>) carrying the new formatting.
    $r = $para.Range
    $r.InsertBefore("X")
    $r.Font.Underline = 1
    $delRng = $d.Range($r.Start, $r.Start + 1)
    $delRng.Delete()
}

# 1. "What type of architecture designs are suitable for the HeardIT application?"
Underline-Paragraph("What type of architecture designs are suitable for the HeardIT application?") | Out-Null

# 2. "Which architecture design meets the scalability and stability requirements?"
Underline-Paragraph("Which architecture design meets the scalability and stability requirements?") | Out-Null

# 3. "Which architecture design meets the maintainability requirements?" plus the
#    blank paragraph right before it.
$para3 = Underline-Paragraph("Which architecture design meets the maintainability requirements?")
if ($para3) {
    $prevPara3 = $para3.Previous()
    Underline-EmptyParagraph($prevPara3)
}

# 4. "Which architecture design meets the requirements for cloud native applications?"
Underline-Paragraph("Which architecture design meets the requirements for cloud native applications?") | Out-Null

# 5. "Which architecture design meets the security requirements?"
Underline-Paragraph("Which architecture design meets the security requirements?") | Out-Null

Write-Host "Underline edits applied"
